# Add extra data to Race line excel sheet
# Adds two new columns (J, K) of integer index/grouping data to the
# "Raceline Data" sheet, rows 4-42, formatted as whole numbers, right
# aligned / vertically centered, and widens column J to fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raceline Data")

# --- New column data (rows 4..42) -----------------------------------
$jVals = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38)
$kVals = @(0,1,2,3,4,4,4,5,6,7,7,7,8,9,10,11,12,12,12,13,14,15,16,16,16,17,18,19,20,20,20,21,22,23,24,25,25,25,26)
$n = $jVals.Length

$jData = New-Object 'object[,]' $n,1
$kData = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) {
    $jData[$i,0] = $jVals[$i]
    $kData[$i,0] = $kVals[$i]
}
$ws.Range("J4:J42").Value = $jData
$ws.Range("K4:K42").Value = $kData

# --- Formatting: whole-number format, right/center aligned -----------
# Build the format on a single cell first, then fan it out with a
# format-only paste so the whole block shares one cell style instead of
# each property-assignment minting its own style entry.
$formatCell = $ws.Cells.Item(4, 10)
$formatCell.NumberFormat = "0"
$formatCell.HorizontalAlignment = -4152   # xlRight
$formatCell.VerticalAlignment = -4108     # xlCenter
$formatCell.Copy() | Out-Null
$ws.Range("J4:K42").PasteSpecial(-4122) | Out-Null      # xlPasteFormats
$excel.CutCopyMode = $false

# --- Column J width (best-fit to the new integer data) ---------------
$ws.Columns.Item(10).AutoFit() | Out-Null

# --- Selection / scroll, matching where the author ended up ----------
$ws.Range("K41").Select()
